$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# New test cases ENWIAM100 / ENWIAM101 are appended as rows 23 and 24,
# matching the bordered / non-filled look of the existing table rows
# (row 21 is a plain, unshaded data row so it makes a good formatting
# template for the two new rows).
$ws.Range("A21:E21").Copy()
$ws.Range("A23:E23").PasteSpecial(-4122)
$ws.Range("A21:E21").Copy()
$ws.Range("A24:E24").PasteSpecial(-4122)

# Column C in row 23 should not wrap (unlike most Description cells), so
# pull its formatting from a non-wrapping column instead.
$ws.Range("A21").Copy()
$ws.Range("C23").PasteSpecial(-4122)

# Fill in the new TCIDs first (ENWIAM100 / ENWIAM101), then complete the
# rest of row 23, then the rest of row 24.
$ws.Range("A23").Value = "ENWIAM100"
$ws.Range("A24").Value = "ENWIAM101"

# Row 23: ENWIAM100
$ws.Range("B23").Value = "OPQA-2119||OPQA-2287||OPQA-2293||OPQA-2305||OPQA-2308||OPQA-2319||OPQA-2336"
$ws.Range("C23").Value = "Sign-in with social and link existing steam account with matching email."
$ws.Range("D23").Value = "Y"

# Row 24: ENWIAM101
$ws.Range("B24").Value = "OBT"
$ws.Range("C24").Value = "OB T"
$ws.Range("D24").Value = "Y"

$null = $ws.Range("C30").Select()
